# "Filters and export fix - missing filtered output"
#
# The filters template had three helper/example columns - "name1",
# "name2" and "description" - that were never actually used by the
# filtering logic and were hiding the real "descTerm" output column.
# Drop them so the table only exposes the four columns that matter:
# searchColumn_contains, searchTerm, validate, descTerm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# "descTerm" (originally column F) is the one column after the three
# helper columns that needs to survive. Move it next to "validate" so
# the three columns being dropped ("name1", "name2", "description") end
# up contiguous (E:G) and can be removed in a single delete.
$ws.Range("F1").EntireColumn.Cut()
$ws.Range("D1").EntireColumn.Insert()

# Remove the now-contiguous "name1", "name2", "description" columns.
$ws.Range("E1:G1").EntireColumn.Delete()

# Re-sync the table/autoFilter definition to the new A1:D2 extent.
$lo.Resize($ws.Range("A1:D2"))

# Re-stamp the surviving header so the table picks up the "descTerm"
# column name instead of the placeholder name it inherited positionally.
$ws.Range("D1").Value = "descTerm"

# Restore the active cell selection.
$ws.Range("F6").Select()
